$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.487.69"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.654.04"
$ws.Range("E3").Value = "  -3.89%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  -1.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.66"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -1.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3632"
$ws.Range("E7").Value = "  -4.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.40"
$ws.Range("E8").Value = "  -4.51%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3287"
$ws.Range("E9").Value = "  -6.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.121"
$ws.Range("E10").Value = "  -5.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06933"
$ws.Range("E11").Value = "  -7.62%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9968"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.964"
$ws.Range("E13").Value = "  -6.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.20"
$ws.Range("E14").Value = "  -7.91%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.649.83"
$ws.Range("E15").Value = "  -4.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.567"
$ws.Range("E16").Value = "  -6.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001048"
$ws.Range("E17").Value = "  -7.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06487"
$ws.Range("E18").Value = "  -3.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9990"
$ws.Range("E19").Value = "  -1.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.08"
$ws.Range("E20").Value = "  -8.77%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.925"
$ws.Range("E21").Value = "  -7.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.80"
$ws.Range("E22").Value = "  -8.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.22"
$ws.Range("E23").Value = "  -6.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.428.19"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.417"
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.346"
$ws.Range("E26").Value = "  -16.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "145.87"
$ws.Range("E27").Value = "  -3.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.51"
$ws.Range("E28").Value = "  -9.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.834.10"
$ws.Range("E29").Value = "  -4.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "124.89"
$ws.Range("E30").Value = "  -5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.162"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.046"
$ws.Range("E32").Value = "  -4.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.612"
$ws.Range("E33").Value = "  -18.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08336"
$ws.Range("E34").Value = "  -4.41%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.668"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.31"
$ws.Range("E36").Value = "  -10.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.209"
$ws.Range("E37").Value = "  -6.99%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06058"
$ws.Range("E38").Value = "  -7.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02216"
$ws.Range("E39").Value = "  -9.37%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.246"
$ws.Range("E40").Value = "  -9.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.205"
$ws.Range("E41").Value = "  -5.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2045"
$ws.Range("E42").Value = "  -7.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9984"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5863"
$ws.Range("E44").Value = "  -9.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.723"
$ws.Range("E45").Value = "  -3.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.68"
$ws.Range("E46").Value = "  -10.00%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5582"
$ws.Range("E47").Value = "  -9.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.06"
$ws.Range("E48").Value = "  -5.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.937"
$ws.Range("E49").Value = "  -10.08%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06904"
$ws.Range("E50").Value = "  -5.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.07"
$ws.Range("E51").Value = "  -7.30%  "
